# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: zoom level + selection -------------------------------
# Before: zoomScale/zoomScaleNormal = 190, selection F2 (F2:F11)
# After:  zoomScale/zoomScaleNormal = 235, selection A1 (A1:G1)
$excel.ActiveWindow.Zoom = 235
$ws.Range("A1:G1").Select()

# --- Column widths ------------------------------------------------------
# Column A: 18.2166666666667 -> 12.625
# Columns F:G used to be merged into one <col min="6" max="7"> entry at
# 13.8416666666667; they now need distinct widths (F=8.125, G=11.625).
# ColumnWidth values are quantized by Excel to whole pixels (character
# width grid, 1/7 increments with the default font), so we feed in the
# character-width value whose rounded pixel width lands closest to the
# target stored width.
$ws.Columns.Item(1).ColumnWidth = 11.857142857142858   # -> stored 12.571428571428571 (closest to 12.625)
$ws.Columns.Item(6).ColumnWidth = 7.428571428571429    # -> stored 8.142857142857142  (closest to 8.125)
$ws.Columns.Item(7).ColumnWidth = 10.857142857142858   # -> stored 11.571428571428571 (closest to 11.625)
